$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Time entries for row 15 (Tuesday) ---
# C15 = 1:30 PM (13.5/24), D15 = 5:30 PM (17.5/24)
$ws.Range("C15").Value = 0.5625
$ws.Range("D15").Value = 0.729166666666667

# --- Selection moves from D15 to D16 ---
$ws.Range("D16").Select()

# --- Column width tweaks (columns narrowed slightly) ---
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 5.666666666666667
$ws.Columns.Item(8).ColumnWidth = 6.666666666666667
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 5.666666666666667

# --- Print area set again twice, appending two more generations of the
#     auto-numbered "_xlnm.Print_Area_0..._0" defined name that LibreOffice/
#     Excel mint every time the print area is (re)applied ---
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")
